$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$row = @("M1", "Tnfsf18", "Tnfrsf18", "ECs", 1, 0.3333333333333333, 0.1977553333333333, 0.593266, 0.1658261216405568, 0.1658261216405568, 2, 1, 0.8688675, 1.737735, 0.2876230881030516, 0.2287602607316837, 0.171823182085, 1.03093909251, 0.04769542119440923, 0.03793442682261768)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(2, $c).Value() = $row[$c-1] }

# Row 3
$row = @("M1", "Tnfsf18", "Tnfrsf18", "FAPs", 1, 0.3333333333333333, 0.1977553333333333, 0.593266, 0.1658261216405568, 0.1658261216405568, 3, 1, 0.524837, 1.574511, 0.1737379274639013, 0.207273000132301, 0.1037893158806667, 0.9341038429259999, 0.02881028669320714, 0.03437127773274209)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(3, $c).Value() = $row[$c-1] }

# Row 4
$row = @("M1", "Tnfsf18", "Tnfrsf18", "M1", 1, 0.3333333333333333, 0.1977553333333333, 0.593266, 0.1658261216405568, 0.1658261216405568, 2, 0.6666666666666666, 0.2588029999999999, 0.7764089999999999, 0.08567211694571848, 0.1022086366876571, 0.05117967353266665, 0.4606170617939999, 0.01420667488584472, 0.01694886182008291)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(4, $c).Value() = $row[$c-1] }

# Row 5
$row = @("M1", "Tnfsf18", "Tnfrsf18", "M2", 1, 0.3333333333333333, 0.1977553333333333, 0.593266, 0.1658261216405568, 0.1658261216405568, 3, 1, 0.574252, 1.722756, 0.1900958818109247, 0.2267883835780902, 0.1135613956773333, 1.022052561096, 0.0315228628205473, 0.03760743808188564)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(5, $c).Value() = $row[$c-1] }

# Row 6
$row = @("M1", "Tnfsf18", "Tnfrsf18", "Neutro", 1, 0.3333333333333333, 0.1977553333333333, 0.593266, 0.1658261216405568, 0.1658261216405568, 3, 1, 0.196714, 0.5901420000000001, 0.06511866096165837, 0.07768793158261605, 0.03890124264133334, 0.350111183772, 0.01079837499369814, 0.01288268839262215)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(6, $c).Value() = $row[$c-1] }

# Row 7
$row = @("M1", "Tnfsf18", "Tnfrsf18", "sCs", 1, 0.3333333333333333, 0.1977553333333333, 0.593266, 0.1658261216405568, 0.1658261216405568, 2, 1, 0.5973809999999999, 1.194762, 0.1977523247147454, 0.157281787287652, 0.118135278782, 0.7088116726919999, 0.03279250105285026, 0.02608142879060636)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(7, $c).Value() = $row[$c-1] }

# Row 8
$row = @("M2", "Tnfsf18", "Tnfrsf18", "ECs", 3, 1, 0.9947909999999999, 2.984373, 0.8341738783594431, 0.8341738783594432, 2, 1, 0.8688675, 1.737735, 0.2876230881030516, 0.2287602607316837, 0.8643415691924999, 5.186049415155, 0.2399276669086424, 0.1908258339090661)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(8, $c).Value() = $row[$c-1] }

# Row 9
$row = @("M2", "Tnfsf18", "Tnfrsf18", "FAPs", 3, 1, 0.9947909999999999, 2.984373, 0.8341738783594431, 0.8341738783594432, 3, 1, 0.524837, 1.574511, 0.1737379274639013, 0.207273000132301, 0.5221031240669999, 4.698928116603, 0.1449276407706942, 0.1729017223995589)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(9, $c).Value() = $row[$c-1] }

# Row 10
$row = @("M2", "Tnfsf18", "Tnfrsf18", "M1", 3, 1, 0.9947909999999999, 2.984373, 0.8341738783594431, 0.8341738783594432, 2, 0.6666666666666666, 0.2588029999999999, 0.7764089999999999, 0.08567211694571848, 0.1022086366876571, 0.2574548951729999, 2.317094056557, 0.07146544205987375, 0.0852597748675742)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(10, $c).Value() = $row[$c-1] }

# Row 11
$row = @("M2", "Tnfsf18", "Tnfrsf18", "M2", 3, 1, 0.9947909999999999, 2.984373, 0.8341738783594431, 0.8341738783594432, 3, 1, 0.574252, 1.722756, 0.1900958818109247, 0.2267883835780902, 0.5712607213319999, 5.141346491987999, 0.1585730189903774, 0.1891809454962045)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(11, $c).Value() = $row[$c-1] }

# Row 12
$row = @("M2", "Tnfsf18", "Tnfrsf18", "Neutro", 3, 1, 0.9947909999999999, 2.984373, 0.8341738783594431, 0.8341738783594432, 3, 1, 0.196714, 0.5901420000000001, 0.06511866096165837, 0.07768793158261605, 0.195689316774, 1.761203850966, 0.05432028596796023, 0.06480524318999391)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(12, $c).Value() = $row[$c-1] }

# Row 13
$row = @("M2", "Tnfsf18", "Tnfrsf18", "sCs", 3, 1, 0.9947909999999999, 2.984373, 0.8341738783594431, 0.8341738783594432, 2, 1, 0.5973809999999999, 1.194762, 0.1977523247147454, 0.157281787287652, 0.5942692423709999, 3.565615454225999, 0.1649598236618952, 0.1312003584970456)
for ($c = 1; $c -le $row.Length; $c++) { $ws.Cells.Item(13, $c).Value() = $row[$c-1] }
